# Weekly update: insert two new daily price records for
# "Vega Modelo de Temuco - Perejil" ahead of the existing block
# (previous rows 358-383 shift down to 360-385), extending the
# used range from A1:R383 to A1:R385.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 358/359, pushing old rows 358-383 down to 360-385.
$ws.Range("A358:A359").EntireRow.Insert()

# --- New row 358 ---
$ws.Cells.Item(358, 1).Value = 10
$ws.Cells.Item(358, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(358, 3).Value = "La Araucanía"
$ws.Cells.Item(358, 4).Value = 44826
$ws.Cells.Item(358, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(358, 5).Value = 9
$ws.Cells.Item(358, 6).Value = 100112044
$ws.Cells.Item(358, 7).Value = "Perejil"
$ws.Cells.Item(358, 8).Value = "Sin especificar"
$ws.Cells.Item(358, 9).Value = "Primera"
$ws.Cells.Item(358, 10).Value = 50
$ws.Cells.Item(358, 11).Value = 4000
$ws.Cells.Item(358, 12).Value = 4000
$ws.Cells.Item(358, 13).Value = 4000
$ws.Cells.Item(358, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(358, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(358, 16).Value = 1333
$ws.Cells.Item(358, 17).Value = 3
$ws.Cells.Item(358, 18).Value = "Hortaliza"

# --- New row 359 ---
$ws.Cells.Item(359, 1).Value = 10
$ws.Cells.Item(359, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(359, 3).Value = "La Araucanía"
$ws.Cells.Item(359, 4).Value = 44826
$ws.Cells.Item(359, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(359, 5).Value = 9
$ws.Cells.Item(359, 6).Value = 100112044
$ws.Cells.Item(359, 7).Value = "Perejil"
$ws.Cells.Item(359, 8).Value = "Sin especificar"
$ws.Cells.Item(359, 9).Value = "Primera"
$ws.Cells.Item(359, 10).Value = 80
$ws.Cells.Item(359, 11).Value = 3300
$ws.Cells.Item(359, 12).Value = 3300
$ws.Cells.Item(359, 13).Value = 3300
$ws.Cells.Item(359, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(359, 15).Value = "Región Metropolitana"
$ws.Cells.Item(359, 16).Value = 1100
$ws.Cells.Item(359, 17).Value = 3
$ws.Cells.Item(359, 18).Value = "Hortaliza"
